$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NPLQ fixtures")

$ws.Range("I109").Value = "Y"

$ws.Range("C110").Value = "Logan Lightning"
$ws.Range("D110").Value = "LIG"
$ws.Range("E110").Value = 3
$ws.Range("F110").Value = "Capalaba FC"
$ws.Range("G110").Value = "CAP"
$ws.Range("H110").Value = 2
$ws.Range("I110").Value = "Y"

$ws.Range("C111").Value = "Eastern Suburbs"
$ws.Range("D111").Value = "EAS"
$ws.Range("E111").Value = 1
$ws.Range("F111").Value = "Brisbane City"
$ws.Range("G111").Value = "BCT"
$ws.Range("H111").Value = 4
$ws.Range("I111").Value = "Y"

$ws.Range("B112").Value = 44780
$ws.Range("C112").Value = "Moreton Bay United"
$ws.Range("D112").Value = "MBJ"
$ws.Range("F112").Value = "Peninsula Power"
$ws.Range("G112").Value = "PEN"
$ws.Range("H112").Value = 2
$ws.Range("I112").Value = "Y"

$ws.Range("E113").Value = 4
$ws.Range("H113").Value = 2
$ws.Range("I113").Value = "Y"

$ws.Range("E114").Value = 2
$ws.Range("H114").Value = 3
$ws.Range("I114").Value = "Y"

$ws.Range("E115").Value = 3
$ws.Range("H115").Value = 1

$ws.Range("E116").Value = 1
$ws.Range("H116").Value = 2

$ws.Range("C117").Value = "Sunshine Coast Wanderers"
$ws.Range("D117").Value = "SCW"
$ws.Range("E117").Value = 0
$ws.Range("F117").Value = "Gold Coast United"
$ws.Range("G117").Value = "GCU"
$ws.Range("H117").Value = 0

$ws.Range("C118").Value = "Brisbane City"
$ws.Range("D118").Value = "BCT"
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = "Capalaba FC"
$ws.Range("G118").Value = "CAP"
$ws.Range("H118").Value = 0

$ws.Range("C119").Value = "Gold Coast Knights"
$ws.Range("D119").Value = "GCK"
$ws.Range("E119").Value = 2
$ws.Range("F119").Value = "Olympic FC"
$ws.Range("G119").Value = "BOL"
$ws.Range("H119").Value = 2

$ws.Range("C120").Value = "Eastern Suburbs"
$ws.Range("D120").Value = "EAS"
$ws.Range("E120").Value = 1
$ws.Range("F120").Value = "Brisbane Roar Youth"
$ws.Range("G120").Value = "BRR"
$ws.Range("H120").Value = 0

$ws.Range("C121").Value = "Lions"
$ws.Range("D121").Value = "LIO"
$ws.Range("E121").Value = 3
$ws.Range("F121").Value = "Moreton Bay United"
$ws.Range("G121").Value = "MBJ"
$ws.Range("H121").Value = 1

$ws.Range("B122").Value = 44790
$ws.Range("C122").Value = "Peninsula Power"
$ws.Range("D122").Value = "PEN"
$ws.Range("F122").Value = "Logan Lightning"
$ws.Range("G122").Value = "LIG"

$ws.Range("C123").Value = "Brisbane Roar Youth"
$ws.Range("D123").Value = "BRR"
$ws.Range("F123").Value = "Gold Coast Knights"
$ws.Range("G123").Value = "GCK"

$ws.Range("A124").Value = 22
$ws.Range("B124").Value = 44793
$ws.Range("C124").Value = "Gold Coast United"
$ws.Range("D124").Value = "GCU"
$ws.Range("F124").Value = "Eastern Suburbs"
$ws.Range("G124").Value = "EAS"

$ws.Range("C125").Value = "Brisbane City"
$ws.Range("D125").Value = "BCT"
$ws.Range("F125").Value = "Brisbane Roar Youth"
$ws.Range("G125").Value = "BRR"

$ws.Range("C126").Value = "Gold Coast Knights"
$ws.Range("D126").Value = "GCK"
$ws.Range("F126").Value = "Moreton Bay United"
$ws.Range("G126").Value = "MBJ"

$ws.Range("C127").Value = "Olympic FC"
$ws.Range("D127").Value = "BOL"
$ws.Range("F127").Value = "Sunshine Coast Wanderers"
$ws.Range("G127").Value = "SCW"

$ws.Range("C128").Value = "Lions"
$ws.Range("D128").Value = "LIO"
$ws.Range("F128").Value = "Logan Lightning"
$ws.Range("G128").Value = "LIG"

$ws.Range("C129").Value = "Capalaba FC"
$ws.Range("D129").Value = "CAP"
$ws.Range("F129").Value = "Peninsula Power"
$ws.Range("G129").Value = "PEN"

$ws.Range("A130").Value = 18
$ws.Range("B130").Value = 44796
$ws.Range("C130").Value = "Gold Coast United"
$ws.Range("D130").Value = "GCU"

$ws.Range("A131").Value = 8
$ws.Range("B131").Value = 44800
$ws.Range("C131").Value = "Logan Lightning"
$ws.Range("D131").Value = "LIG"
$ws.Range("F131").Value = "Brisbane City"
$ws.Range("G131").Value = "BCT"

$ws.Range("B116").Select()

Write-Output "Done applying NPLQ fixtures update"